$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 319, shifting existing rows 319:423 down to 321:425
$ws.Rows("319:320").Insert()

# Row 319 - new data row (Primera)
$ws.Range("A319").Value = 8
$ws.Range("B319").Value = "Terminal La Palmera de La Serena"
$ws.Range("C319").Value = "Coquimbo"
$ws.Range("D319").Value = 44722
$ws.Range("E319").Value = 4
$ws.Range("F319").Value = 100112043
$ws.Range("G319").Value = "Pepino dulce"
$ws.Range("H319").Value = "Cultivar IV Región"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 500
$ws.Range("K319").Value = 13000
$ws.Range("L319").Value = 14000
$ws.Range("M319").Value = 13500
$ws.Range("N319").Value = '$/bandeja 18 kilos'
$ws.Range("O319").Value = "Provincia de Limarí"
$ws.Range("P319").Value = 750
$ws.Range("Q319").Value = 18
$ws.Range("R319").Value = "Hortaliza"

# Row 320 - new data row (Segunda)
$ws.Range("A320").Value = 8
$ws.Range("B320").Value = "Terminal La Palmera de La Serena"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 44722
$ws.Range("E320").Value = 4
$ws.Range("F320").Value = 100112043
$ws.Range("G320").Value = "Pepino dulce"
$ws.Range("H320").Value = "Cultivar IV Región"
$ws.Range("I320").Value = "Segunda"
$ws.Range("J320").Value = 300
$ws.Range("K320").Value = 10000
$ws.Range("L320").Value = 11000
$ws.Range("M320").Value = 10500
$ws.Range("N320").Value = '$/bandeja 18 kilos'
$ws.Range("O320").Value = "Provincia de Limarí"
$ws.Range("P320").Value = 583
$ws.Range("Q320").Value = 18
$ws.Range("R320").Value = "Hortaliza"
